$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.265.47"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -4.54%  "

$ws.Range("D3").Value = "'1.855.77"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -5.73%  "

$ws.Range("D4").Value = "'1.001"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.97%  "

$ws.Range("D5").Value = "'321.62"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.61%  "

$ws.Range("E6").Value = "  -0.99%  "

$ws.Range("D7").Value = "'0.4493"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -6.39%  "

$ws.Range("D8").Value = "'0.3846"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -5.51%  "

$ws.Range("D9").Value = "'47.94"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -11.40%  "

$ws.Range("D10").Value = "'0.07874"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -7.59%  "

$ws.Range("E11").Value = "  -4.52%  "

$ws.Range("D12").Value = "'21.29"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -5.46%  "

$ws.Range("D13").Value = "'1.869.44"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -6.33%  "

$ws.Range("D14").Value = "'5.866"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -5.55%  "

$ws.Range("D15").Value = "'7.144"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -6.65%  "

$ws.Range("D16").Value = "'0.9999"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.28%  "

$ws.Range("E17").Value = "  -4.30%  "

$ws.Range("D18").Value = "'85.38"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -6.60%  "

$ws.Range("D19").Value = "'0.06527"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.84%  "

$ws.Range("D20").Value = "'16.89"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -9.56%  "

$ws.Range("D21").Value = "'1.001"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.03%  "

$ws.Range("D22").Value = "'5.501"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -6.52%  "

$ws.Range("D23").Value = "'27.288.92"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -4.61%  "

$ws.Range("D24").Value = "'10.73"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -6.96%  "

$ws.Range("D25").Value = "'2.262"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.65%  "

$ws.Range("D26").Value = "'2.077.52"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -6.66%  "

$ws.Range("D27").Value = "'151.77"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.02%  "

$ws.Range("D28").Value = "'19.67"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.65%  "

$ws.Range("E29").Value = "  -6.16%  "

$ws.Range("D30").Value = "'5.428"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -8.22%  "

$ws.Range("D31").Value = "'120.19"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.89%  "

$ws.Range("B32").Value = "ARBITRUM"
$ws.Range("C32").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D32").Value = "'1.470"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.25%  "

$ws.Range("B33").Value = "Stellar"
$ws.Range("C33").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D33").Value = "'0.09265"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -4.34%  "

$ws.Range("B34").Value = "ImmutableX"
$ws.Range("C34").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D34").Value = "'0.9350"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -5.72%  "

$ws.Range("D35").Value = "'3.582"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.14%  "

$ws.Range("D36").Value = "'5.276"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -6.64%  "

$ws.Range("E37").Value = "  -5.17%  "

$ws.Range("D38").Value = "'0.05975"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -4.50%  "

$ws.Range("D39").Value = "'1.209"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -4.11%  "

$ws.Range("D40").Value = "'8.248"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -10.00%  "

$ws.Range("D41").Value = "'0.9999"
$ws.Range("D41").Style = "Normal"

$ws.Range("D42").Value = "'0.5902"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -5.51%  "

$ws.Range("D43").Value = "'0.1875"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.66%  "

$ws.Range("E44").Value = "  -10.37%  "

$ws.Range("D45").Value = "'1.260"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -6.95%  "

$ws.Range("D46").Value = "'0.5623"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -5.77%  "

$ws.Range("E47").Value = "  -9.13%  "

$ws.Range("D48").Value = "'3.354"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.91%  "

$ws.Range("E49").Value = "  -7.80%  "

$ws.Range("D50").Value = "'0.06799"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.38%  "

$ws.Range("D51").Value = "'107.99"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.22%  "
